# Update column F ("dSF") values per repulled/recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    3  = 8
    4  = -6
    5  = -1
    6  = -2
    7  = 1
    8  = 3
    9  = -1
    10 = 1
    11 = 2
    12 = 2
    13 = 1
    14 = -2
    15 = -3
    16 = -1
    18 = 1
    19 = -1
    20 = -2
    21 = -2
    22 = 8
    23 = 8
    25 = 6
    26 = 2
    27 = 0
    28 = -2
    29 = 1
    32 = -1
    33 = 1
    34 = -1
    35 = 3
    36 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
